$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the style of the "free-form" columns (A:E study/n/N, I:N corrected..clinical_setting)
# so that, once their contents are cleared, the cells disappear entirely (default style),
# matching rows whose data was removed for the new extraction round.
$ws.Range("A30:E32").Style = "Normal"
$ws.Range("I30:N32").Style = "Normal"

# Clear all the values for the old studies (Sastre 2019 / Harioka 2000) on rows 30-32,
# across every column (A:T) - the RoB/funding columns (F:H, O:T) keep their formatting
# since they still carry a non-default style.
$ws.Range("A30:T32").ClearContents()

# Row height was an explicit 16pt for the populated rows; with the rows now blank,
# restore the sheet's default row height.
$ws.Rows("30:32").AutoFit()

# Move the selection to the next empty row, ready for the new "R1" extraction entries.
$ws.Range("A33").Select()
